# Apply cryptos list update - generated from diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '30.717.69'
$ws.Range('E2').Value2 = '  +2.29%  '
$ws.Range('D3').Value2 = '1.890.79'
$ws.Range('E3').Value2 = '  +0.79%  '
$ws.Range('D4').Value2 = "'1.001"
$ws.Range('E4').Value2 = '  +0.17%  '
$ws.Range('D5').Value2 = "'244.82"
$ws.Range('E5').Value2 = '  +0.67%  '
$ws.Range('D6').Value2 = "'0.9997"
$ws.Range('E6').Value2 = '  +0.09%  '
$ws.Range('D7').Value2 = "'0.4920"
$ws.Range('E7').Value2 = '  +0.18%  '
$ws.Range('D8').Value2 = "'0.2960"
$ws.Range('E8').Value2 = '  +1.15%  '
$ws.Range('D9').Value2 = "'0.06789"
$ws.Range('E9').Value2 = '  +2.91%  '
$ws.Range('D10').Value2 = '1.888.05'
$ws.Range('E10').Value2 = '  +0.66%  '
$ws.Range('D11').Value2 = "'17.28"
$ws.Range('E11').Value2 = '  +3.93%  '
$ws.Range('D12').Value2 = "'0.07227"
$ws.Range('E12').Value2 = '  +0.70%  '
$ws.Range('D13').Value2 = "'91.33"
$ws.Range('E13').Value2 = '  +5.95%  '
$ws.Range('D14').Value2 = "'0.6784"
$ws.Range('E14').Value2 = '  +1.65%  '
$ws.Range('D15').Value2 = "'5.050"
$ws.Range('E15').Value2 = '  +2.77%  '
$ws.Range('D16').Value2 = '30.666.30'
$ws.Range('E16').Value2 = '  +2.34%  '
$ws.Range('D17').Value2 = "'0.000007999"
$ws.Range('E17').Value2 = '  +2.64%  '
$ws.Range('D18').Value2 = "'1.000"
$ws.Range('E18').Value2 = '  +0.12%  '
$ws.Range('D19').Value2 = "'13.16"
$ws.Range('E19').Value2 = '  +2.99%  '
$ws.Range('D20').Value2 = '2.131.56'
$ws.Range('E20').Value2 = '  +0.53%  '
$ws.Range('E21').Value2 = '  +0.29%  '
$ws.Range('D22').Value2 = "'4.827"
$ws.Range('E22').Value2 = '  +1.19%  '
$ws.Range('D23').Value2 = "'190.73"
$ws.Range('E23').Value2 = '  +33.03%  '
$ws.Range('D24').Value2 = "'6.122"
$ws.Range('E24').Value2 = '  +4.40%  '
$ws.Range('D25').Value2 = "'9.373"
$ws.Range('E25').Value2 = '  +3.09%  '
$ws.Range('D26').Value2 = "'155.00"
$ws.Range('E26').Value2 = '  +2.16%  '
$ws.Range('D27').Value2 = "'19.15"
$ws.Range('E27').Value2 = '  +13.17%  '
$ws.Range('D28').Value2 = "'1.907"
$ws.Range('E28').Value2 = '  +0.70%  '
$ws.Range('D29').Value2 = "'1.402"
$ws.Range('E29').Value2 = '  +1.37%  '
$ws.Range('D30').Value2 = "'4.345"
$ws.Range('E30').Value2 = '  +3.61%  '
$ws.Range('D31').Value2 = "'0.09104"
$ws.Range('E31').Value2 = '  +4.11%  '
$ws.Range('D32').Value2 = "'4.017"
$ws.Range('E32').Value2 = '  +1.05%  '
$ws.Range('D33').Value2 = "'0.05210"
$ws.Range('E33').Value2 = '  +4.03%  '
$ws.Range('D34').Value2 = "'0.7575"
$ws.Range('E34').Value2 = '  +5.33%  '
$ws.Range('E35').Value2 = '  +0.07%  '
$ws.Range('D36').Value2 = "'2.773"
$ws.Range('E36').Value2 = '  +4.18%  '
$ws.Range('D37').Value2 = "'0.01844"
$ws.Range('E37').Value2 = '  +1.34%  '
$ws.Range('D38').Value2 = "'2.679"
$ws.Range('E38').Value2 = '  -0.21%  '
$ws.Range('D39').Value2 = "'2.152"
$ws.Range('E39').Value2 = '  -0.12%  '
$ws.Range('D40').Value2 = "'0.9373"
$ws.Range('D41').Value2 = "'0.4428"
$ws.Range('E41').Value2 = '  +5.07%  '
$ws.Range('D42').Value2 = "'105.43"
$ws.Range('E42').Value2 = '  +2.35%  '
$ws.Range('E43').Value2 = '  +0.18%  '
$ws.Range('E44').Value2 = '  +0.17%  '
$ws.Range('D45').Value2 = "'7.619"
$ws.Range('E45').Value2 = '  +3.45%  '
$ws.Range('E46').Value2 = '  +6.08%  '
$ws.Range('D47').Value2 = "'0.05867"
$ws.Range('E47').Value2 = '  +2.96%  '
$ws.Range('D48').Value2 = "'8.762"
$ws.Range('B49').Value2 = 'Decentraland'
$ws.Range('C49').Value2 = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').Value2 = "'0.3945"
$ws.Range('E49').Value2 = '  +4.93%  '
$ws.Range('B50').Value2 = 'NEARProtocol'
$ws.Range('C50').Value2 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value2 = "'1.423"
$ws.Range('E50').Value2 = '  +6.30%  '
$ws.Range('D51').Value2 = "'33.70"
$ws.Range('E51').Value2 = '  +2.91%  '
